# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets
# to reflect refreshed scrape numbers (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6747
$ws1.Range("F3").Value = 12
$ws1.Range("F5").Value = 69
$ws1.Range("F6").Value = 12
$ws1.Range("F8").Value = 99
$ws1.Range("F9").Value = 94
$ws1.Range("F12").Value = 173
$ws1.Range("F13").Value = 408
$ws1.Range("F15").Value = 1613
$ws1.Range("F16").Value = 22
$ws1.Range("F17").Value = 3391
$ws1.Range("F20").Value = 7
$ws1.Range("F21").Value = 2019
$ws1.Range("F22").Value = 134
$ws1.Range("F27").Value = 3

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6747
$ws4.Range("F3").Value = 12
$ws4.Range("F5").Value = 69
$ws4.Range("F6").Value = 12
$ws4.Range("F9").Value = 99
$ws4.Range("F10").Value = 94
$ws4.Range("F13").Value = 173
$ws4.Range("F14").Value = 408
$ws4.Range("F16").Value = 1613
$ws4.Range("F17").Value = 22
$ws4.Range("F18").Value = 3391
$ws4.Range("F21").Value = 7
$ws4.Range("F22").Value = 2019
$ws4.Range("F23").Value = 134
$ws4.Range("F28").Value = 3
